# Trade #17 closed at 2026-02-17 04:07:52 - unknown UNKNOWN +0.000%
#
# Adds the new closed trade (Trade #17) to the "All Trades" and
# "MarketMaking" sheets, and rolls the updated aggregate stats into the
# "Summary" and "Strategy Status" sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.01   # Current Capital
$summary.Range("B4").Value = 0.01      # Total P&L $
$summary.Range("B5").Value = 0.01      # Total P&L %
$summary.Range("B6").Value = 17        # Total Trades
$summary.Range("B8").Value = 7         # Losing Trades
$summary.Range("B9").Value = 29.41     # Win Rate %

# ---------------------------------------------------------------------------
# Strategy Status sheet (MarketMaking row)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.01     # Capital
$status.Range("D4").Value = 17         # Trades
$status.Range("E4").Value = 0.01       # P&L $
$status.Range("F4").Value = 0.01       # P&L %
$status.Range("G4").Value = 29.41      # Win Rate %

# ---------------------------------------------------------------------------
# New trade row (#17) appended to "All Trades" and "MarketMaking" sheets
# ---------------------------------------------------------------------------
$tradeSheets = @("All Trades", "MarketMaking")

foreach ($sheetName in $tradeSheets) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Force column B to stay plain text (matches existing rows) instead of
    # being auto-parsed into a date serial by the "2026-02-17" literal.
    $ws.Cells.Item(18, 2).NumberFormat = "@"

    $ws.Cells.Item(18, 1).Value = 17
    $ws.Cells.Item(18, 2).Value = "2026-02-17"
    $ws.Cells.Item(18, 3).Value = "04:07:47"
    $ws.Cells.Item(18, 4).Value = "MarketMaking"
    $ws.Cells.Item(18, 5).Value = "UP"
    $ws.Cells.Item(18, 6).Value = 0.19
    $ws.Cells.Item(18, 7).Value = 0.18
    $ws.Cells.Item(18, 8).Value = "CLOSED"
    $ws.Cells.Item(18, 9).Value = -5.2632
    $ws.Cells.Item(18, 10).Value = -0.01
    $ws.Cells.Item(18, 11).Value = 100.01
    $ws.Cells.Item(18, 12).Value = 0
    $ws.Cells.Item(18, 13).Value = 0
    $ws.Cells.Item(18, 14).Value = 0.6
    $ws.Cells.Item(18, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item(18, 16).Value = "early_exit"
    $ws.Cells.Item(18, 17).Value = 0.11
}
